$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$excel.ActiveWindow.Zoom = 85

$ws.Range("J32").Value = 3
$ws.Range("K32").Value = 2
$ws.Range("L32").Value = 3

$ws.Range("J33").Value = 2
$ws.Range("K33").Value = 3
$ws.Range("L33").Value = 2

$ws.Range("J35").Value = 2
$ws.Range("K35").Value = 3
$ws.Range("L35").Value = 3

$ws.Range("L38").Select()
